$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6351431704.04
$ws.Range("C3").Value = 6645541.1728499997

$ws.Range("F17").Select()
